$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace "image" with "video" in D5 (TestCase_004 row, media_type column)
$ws.Range("D5").Value = "video"

# Update the active selection to D6 (mirrors the saved selection state in the diff)
$ws.Range("D6").Select()
